# Applies the "cicada" regeneration edit to the begin-age StructureDefinition
# workbook: updates the URL + Date metadata values, inserts a new
# "Jurisdiction" metadata row, and keeps the Elements sheet's Fixed Value in
# sync (it shares the same URL text).

$wb = $excel.ActiveWorkbook

$metaSheet = $wb.Worksheets.Item("Metadata")
$elemSheet = $wb.Worksheets.Item("Elements")

# Insert a new row above row 11 ("Description") so the metadata table gains a
# "Jurisdiction" entry between "Contact" and "Description".
$metaSheet.Rows.Item(11).Insert()

$metaSheet.Cells.Item(11, 1).Value = "Jurisdiction"
$metaSheet.Cells.Item(11, 2).Value = ""

# Update the URL (pythia -> cicada) and regeneration Date.
$metaSheet.Cells.Item(2, 2).Value = "http://fhirfli.dev/fhir/ig/cicada/StructureDefinition/begin-age"
$metaSheet.Cells.Item(8, 2).Value = "2026-02-11T14:37:07-05:00"

# The Elements sheet's "Fixed Value" column for the Extension.url row holds
# the same URL text; keep it consistent with the Metadata sheet.
$elemSheet.Cells.Item(5, 18).Value = "http://fhirfli.dev/fhir/ig/cicada/StructureDefinition/begin-age"
